$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 to be a text value "abc123" instead of the numeric 1234
$ws.Range("B2").Value = "abc123"

# Add a new row 3 with a second username/password pair
$ws.Range("A3").Value = "nimal"
$ws.Range("B3").Value = "nim123"

# Select B3 as the active cell, matching the resulting selection state
$ws.Range("B3").Select()
